$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Number of Enrollees" - bump the 2024-2nd-semester total (new enrollee)
# ---------------------------------------------------------------------------
$wsEnrollees = $wb.Worksheets.Item("Number of Enrollees")
$wsEnrollees.Range("B5").Value = 5

# ---------------------------------------------------------------------------
# 2) "Number of Enrollees Per Program" - BSCS count goes up by one
# ---------------------------------------------------------------------------
$wsProgram = $wb.Worksheets.Item("Number of Enrollees Per Program")
$wsProgram.Range("B2").Value = 4

# ---------------------------------------------------------------------------
# 3) "Per College" - COS count goes up by one
# ---------------------------------------------------------------------------
$wsCollege = $wb.Worksheets.Item("Per College")
$wsCollege.Range("B2").Value = 5

# ---------------------------------------------------------------------------
# 4) "Per Period" - 2024 2nd semester / COS / BSCS count goes up by one
# ---------------------------------------------------------------------------
$wsPeriod = $wb.Worksheets.Item("Per Period")
$wsPeriod.Range("B50").Value = 2

# ---------------------------------------------------------------------------
# 5) "Age" - a new student (age 22) is added for 2024-2nd; insert a row so the
#    data stays sorted, shifting the existing "2024-2nd-23" row down.
# ---------------------------------------------------------------------------
$wsAge = $wb.Worksheets.Item("Age")
$wsAge.Range("A5").EntireRow.Insert()
$wsAge.Range("A5").Value = "2024-2nd-22"
$wsAge.Range("B5").Value = 1

# ---------------------------------------------------------------------------
# 6) "Gender" - Female count for 2024-2nd semester goes up by one
# ---------------------------------------------------------------------------
$wsGender = $wb.Worksheets.Item("Gender")
$wsGender.Range("C5").Value = 2

# ---------------------------------------------------------------------------
# 7) "Faculty - Subjects" - new faculty member appended at the bottom
# ---------------------------------------------------------------------------
$wsFaculty = $wb.Worksheets.Item("Faculty - Subjects")
$wsFaculty.Range("A7").Value = "Dolores Montesines "
$wsFaculty.Range("B7").Value = ""
$wsFaculty.Range("C7").Value = 0
$wsFaculty.Range("D7").Value = 0
$wsFaculty.Range("E7").Value = 0
$wsFaculty.Range("F7").Value = 0

# ---------------------------------------------------------------------------
# 8) "Student - Subjects" - new student appended at the bottom
# ---------------------------------------------------------------------------
$wsStudent = $wb.Worksheets.Item("Student - Subjects")
$wsStudent.Range("A9").Value = "Cruz KC "
$wsStudent.Range("B9").Value = 0
$wsStudent.Range("C9").Value = 0
$wsStudent.Range("D9").Value = 0
$wsStudent.Range("E9").Value = 0

# ---------------------------------------------------------------------------
# 9) "Local or Foreign" - Local count for 2024-2nd semester goes up by one
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("Local or Foreign")
$wsLocal.Range("B5").Value = 2

# ---------------------------------------------------------------------------
# 10) Update the charts whose source ranges grew by one row (Age,
#     Faculty - Subjects, Student - Subjects) so the plotted series cover the
#     newly-added rows.
# ---------------------------------------------------------------------------

# Age chart: one series, columns A (cat) / B (val), now rows 2-6
$ageChart = $wsAge.ChartObjects(1).Chart
$ageChart.SeriesCollection(1).Formula = '=SERIES(Age!$B$1,Age!$A$2:$A$6,Age!$B$2:$B$6,1)'

# Faculty - Subjects chart: series for columns B,C,D,E,F, now rows 2-7
$facultyChart = $wsFaculty.ChartObjects(1).Chart
$facultyChart.SeriesCollection(1).Formula = "=SERIES('Faculty - Subjects'!`$B`$1,'Faculty - Subjects'!`$A`$2:`$A`$7,'Faculty - Subjects'!`$B`$2:`$B`$7,1)"
$facultyChart.SeriesCollection(2).Formula = "=SERIES('Faculty - Subjects'!`$C`$1,'Faculty - Subjects'!`$A`$2:`$A`$7,'Faculty - Subjects'!`$C`$2:`$C`$7,2)"
$facultyChart.SeriesCollection(3).Formula = "=SERIES('Faculty - Subjects'!`$D`$1,'Faculty - Subjects'!`$A`$2:`$A`$7,'Faculty - Subjects'!`$D`$2:`$D`$7,3)"
$facultyChart.SeriesCollection(4).Formula = "=SERIES('Faculty - Subjects'!`$E`$1,'Faculty - Subjects'!`$A`$2:`$A`$7,'Faculty - Subjects'!`$E`$2:`$E`$7,4)"
$facultyChart.SeriesCollection(5).Formula = "=SERIES('Faculty - Subjects'!`$F`$1,'Faculty - Subjects'!`$A`$2:`$A`$7,'Faculty - Subjects'!`$F`$2:`$F`$7,5)"

# Student - Subjects chart: series for columns B,C,D,E, now rows 2-9
$studentChart = $wsStudent.ChartObjects(1).Chart
$studentChart.SeriesCollection(1).Formula = "=SERIES('Student - Subjects'!`$B`$1,'Student - Subjects'!`$A`$2:`$A`$9,'Student - Subjects'!`$B`$2:`$B`$9,1)"
$studentChart.SeriesCollection(2).Formula = "=SERIES('Student - Subjects'!`$C`$1,'Student - Subjects'!`$A`$2:`$A`$9,'Student - Subjects'!`$C`$2:`$C`$9,2)"
$studentChart.SeriesCollection(3).Formula = "=SERIES('Student - Subjects'!`$D`$1,'Student - Subjects'!`$A`$2:`$A`$9,'Student - Subjects'!`$D`$2:`$D`$9,3)"
$studentChart.SeriesCollection(4).Formula = "=SERIES('Student - Subjects'!`$E`$1,'Student - Subjects'!`$A`$2:`$A`$9,'Student - Subjects'!`$E`$2:`$E`$9,4)"
